$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "pre war state added": reset every member's Saves (col B) and Save_Score
# (col C) back to 0 for data rows 2-99, filling in the previously-blank B
# cells along the way.
for ($row = 2; $row -le 99; $row++) {
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
}

# Reflect the saved selection state (cell F19 active) from the workbook.
$ws.Range("F19").Select() | Out-Null
